$p = $ppt.ActivePresentation
$d = $p.Designs.Add("Office Theme", 2)
Write-Host "done"
